$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.501293420791626
$ws.Range("B1").Value = 3.101760625839233
$ws.Range("C1").Value = 4.81257152557373
$ws.Range("D1").Value = 1.921651482582092
$ws.Range("E1").Value = 1.150917410850525
